$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("2025-09-13")
$src.Copy($null, $src)
$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "2025-09-14"

$data = @(
    ,@('いとこのこ', 'いぬちく(著者)', '第38話')
    ,@('宇崎ちゃんは遊びたい！', '丈(著者)', '第127話')
    ,@('不純な彼女達は懺悔しない', 'ポロロッカ(著者)', '第31話')
    ,@('リビルドワールド', '綾村切人(漫画) ナフセ(原作) 吟(キャラクターデザイン) わいっしゅ(世界観デザイン) cell(メカニックデザイン)', '第73話①')
    ,@('アイドル辞めるけど結婚してくれますか!?', '三吉汐美(著者)', '第17話後半')
    ,@('サーシャちゃんとクラスメイトオタクくん', 'はぐはぐ(著者)', '第85話')
    ,@('三枝さんはメガネ先輩と恋を描く', 'セレビィ量産型(著者)', '第22話前半')
    ,@('時間停止勇者―余命３日の設定じゃ世界を救うには短すぎる―', '光永康則', '第６９話『岩鬼停止』②')
    ,@('不徳のギルド', '河添太一', '第９８話：コインの裏表')
    ,@('魔都精兵のスレイブ', '原作:タカヒロ　漫画:竹村洋平', '第161話　神の寵愛')
    ,@('聖液鍛冶屋のエロランタ', 'しげきっくす(著者)', '第8話前半')
    ,@('美人女上司滝沢さん', 'やんBARU(著者)', '第203話')
    ,@('地元のいじめっ子達に仕返ししようとしたら、別の戦いが始まった。', 'マツモトケンゴ', '第64話 更衣室の戦いが始まった（２）')
    ,@('蜘蛛ですが、なにか？', 'かかし朝浩(著者) 馬場翁(原作) 輝竜司(キャラクター原案)', '第76話その1')
    ,@('追放されたチート付与魔術師は 気ままなセカンドライフを謳歌する。', '六志麻あさ 業務用餅 kisui', '第７１話')
    ,@('実は俺、最強でした？', '原作：澄守 彩 漫画：高橋 愛', '第125話　四騎戦の打ち上げです！')
    ,@('勇者に全部奪われた俺は勇者の母親とパーティを組みました！', '久遠まこと(著者) 石のやっさん(原作)', '第30話')
    ,@('アラサーがVTuberになった話。', '犬威赤彦(漫画) とくめい(原作) カラスBTK(キャラクター原案)', '第28話')
    ,@('今日から僕は、彼女の✕✕を解消する。', 'コアヤアコ(著者)', '第3話後半')
    ,@('俺堕ちスレイブヒーローコレクション', 'ゆっ栗栖(著者)', '第13話前半')
    ,@('魔法少女の使い魔なのに今日も怪人が誘惑してくる', 'ぼぶ(著者)', '第2話')
    ,@('解雇された暗黒兵士(30代)のスローなセカンドライフ', '岡沢六十四 るれくちぇ sage・ジョー', '第72話(後編) ダリエルの照らした道')
    ,@('異世界のんびり農家', '剣康之(作画) 内藤騎之介(原作) やすも(キャラクター原案)', '第304話')
    ,@('よくわからないけれど異世界に転生していたようです', '内々けやき あし カオミン', '第139話 よくわからないけれどズゥゥゥンときたみたいです（１）')
    ,@('花子さん、食べないで', '茸谷きの子(著者)', '第10話後半')
    ,@('くらいあの子としたいこと', '碇マナツ(著者)', '第84話')
    ,@('豚のレバーは加熱しろ', 'みなみ(漫画) 逆井卓馬(原作) 遠坂あさぎ(キャラクターデザイン)', '第44話①')
    ,@('勇者のクズ', 'ナカシマ723', '第48話　勇者の矜持 II（後半）')
    ,@('治癒魔法の間違った使い方 ~戦場を駆ける回復要員~', '九我山レキ(漫画) くろかた(原作) ＫｅＧ(キャラクター原案)', '第82話その3')
    ,@('十歳の最強魔導師', '猫月 天乃聖樹', '第4話（後編）')
    ,@('理想のヒモ生活', '日月ネコ(漫画) 渡辺恒彦（ヒーロー文庫／イマジカインフォス）(原作) 文倉十(キャラクター原案)', '第87話　その2')
    ,@('王様ランキング200話～', '十日草輔（とおかそうすけ）', '第267話')
    ,@('ラスボスたちの隠し仔　～魔王城に転生した元社畜プログラマーは自由気ままに『魔導言語《マジックコード》』を開発する～', '森清士郎(漫画) 熊乃げん骨(原作)', '第16話①')
    ,@('異世界転移で女神様から祝福を！～いえ、手持ちの異能があるので結構です～@COMIC', 'コーダ 壁アキオ', '第4話「ドラゴンテイムと異能開眼②」')
    ,@('遊遊じてき。', 'カンケー', '第3話')
    ,@('魔術師クノンは見えている', 'La-na(作画) 南野海風(原作) Ｌａｒｕｈａ(キャラクター原案)', '第40話②')
    ,@('本物のカノジョにしたくなるまで、私で試していいよ。', 'wano(漫画) 有丈ほえる(原作) 緋月ひぐれ(キャラクター原案)', '第10話①')
    ,@('最強勇者パーティーは愛が知りたい', '山田肌襦袢', '第31話「子供に見せれるほうがいい」')
    ,@('リモート授業になったらクラス１の美少女と同居することになった', '原作：三萩せんや(GA文庫／SBクリエイティブ刊)　漫画：jaguchi　キャラクター原案：さとうぽて', '第13話')
    ,@('神への餌は彼女達', '蔵人幸明(原作) まりお(作画)', '第8話前半')
    ,@('辺境モブ貴族のウチに嫁いできた悪役令嬢が、めちゃくちゃできる良い嫁なんだが？', 'tera(原作) 朝倉はやて(作画) 徹田(キャラクター原案)', '番外編①')
    ,@('最強で最速の無限レベルアップ ～スキル【経験値1000倍】と【レベルフリー】でレベル上限の枷が外れた俺は無双する～', 'シオヤマ琴 鳥羽田 航 トモゼロ', '第80話 お見知り推しを！')
    ,@('「門番やってろ」と言われ15年、突っ立ってる間に俺の魔力が9999（最強）に育ってました', '漫画：はり太郎 原作：まさキチ キャラクター原案：カラスBTK', '第3話')
    ,@('終末ツーリング', 'さいとー栄(著者)', '第50話 三沢基地　その６②')
    ,@('ゲーセン少女と異文化交流', '安原宏和(著者)', '第130話')
    ,@('顔に出ない柏田さんと顔に出る太田君＋', '東ふゆ(著者)', '第34話 愛田さんと長津田君の気持ち')
    ,@('修羅幼女の英雄譚～半端者と言われた傭兵、幼女に転生して成り上がる～', '作画：むらたん 原作：沙城流', '第10話(1)')
    ,@('パワハラギルマスをぶん殴ってブラック聖剣ギルドをクビになったので、辺境で聖剣工房を開くことにした', 'だいたいねむい(原作) まお(漫画)', '第10話④')
    ,@('センパイ、自宅警備員の雇用はいかがですか？', '漫画：コブラサナギ 原作：二上圭 キャラ原案：日向あずり', '番外編')
    ,@('姫騎士は蛮族の嫁', 'コトバノリアキ', '第54話① 只人は待望の主')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i + 1
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
    $ws.Cells.Item($row, 4).Value = $data[$i][2]
}

Write-Host "Done. Sheet count:" $wb.Worksheets.Count
